$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("fields")
$ws4 = $wb.Worksheets.Item("options")

# --- sheet "fields": append 3 new field rows (60-62) ---
$ws3.Range("A60").Value = "Giris_Unitesi_CMKS"
$ws3.Range("B60").Value = "Giris_Unitesi_Tipi_CMKS"
$ws3.Range("C60").Value = "Giriş Ünitesi"
$ws3.Range("D60").Value = "select"
$ws3.Range("E60").Value = $true
$ws3.Range("F60").Value = "Giris_Unitesi_Tipi_CMKS_opts"
$ws3.Range("Q60").Value = "radio"

$ws3.Range("A61").Value = "Giris_Unitesi_CMKS"
$ws3.Range("B61").Value = "Sac_Giris_Sekli_CMKS"
$ws3.Range("C61").Value = "Sac Giriş Şekli"
$ws3.Range("D61").Value = "multiselect"
$ws3.Range("E61").Value = $true
$ws3.Range("F61").Value = "Sac_Giris_Sekli_CMKS_opts"
$ws3.Range("Q61").Value = "checkboxes"

$ws3.Range("A62").Value = "Giris_Unitesi_CMKS"
$ws3.Range("B62").Value = "Giris_Unitesi_Motor"
$ws3.Range("C62").Value = "Giriş Ünitesi Motorlu Mu?"
$ws3.Range("D62").Value = "select"
$ws3.Range("E62").Value = $true
$ws3.Range("F62").Value = "Giris_Unitesi_Motor_opts"
$ws3.Range("Q62").Value = "radio"

# --- sheet "options": append 10 new option rows (173-182) ---
# Column B first (top to bottom), then column A, then column C, then column D
# (mirrors the original authoring order so shared-string indices line up)
$ws4.Range("A173").Value = "Giris_Unitesi_Tipi_CMKS_opts"
$ws4.Range("A174").Value = "Giris_Unitesi_Tipi_CMKS_opts"
$ws4.Range("A175").Value = "Giris_Unitesi_Tipi_CMKS_opts"
$ws4.Range("A176").Value = "Giris_Unitesi_Tipi_CMKS_opts"
$ws4.Range("A177").Value = "Giris_Unitesi_Tipi_CMKS_opts"
$ws4.Range("A178").Value = "Giris_Unitesi_Tipi_CMKS_opts"
$ws4.Range("A179").Value = "Sac_Giris_Sekli_CMKS_opts"
$ws4.Range("A180").Value = "Sac_Giris_Sekli_CMKS_opts"
$ws4.Range("A181").Value = "Giris_Unitesi_Motor_opts"
$ws4.Range("A182").Value = "Giris_Unitesi_Motor_opts"

$ws4.Range("B173").Value = "GU63"
$ws4.Range("B174").Value = "GU80"
$ws4.Range("B175").Value = "GU100"
$ws4.Range("B176").Value = "GU125"
$ws4.Range("B177").Value = "GU160"
$ws4.Range("B178").Value = "GUKR140"
$ws4.Range("B179").Value = "A"
$ws4.Range("B180").Value = "U"
$ws4.Range("B181").Value = "Yok"
$ws4.Range("B182").Value = "M"

$ws4.Range("C173").Value = "Ø63 giriş ünitesi"
$ws4.Range("C174").Value = "Ø80 giriş ünitesi"
$ws4.Range("C175").Value = "Ø100 giriş ünitesi"
$ws4.Range("C176").Value = "Ø125 giriş ünitesi"
$ws4.Range("C177").Value = "Ø160 giriş ünitesi"
$ws4.Range("C178").Value = "Ø140 bilmem neli giriş ünitesi"
$ws4.Range("C179").Value = "Alttan Besleme"
$ws4.Range("C180").Value = "Üstten Besleme"
$ws4.Range("C181").Value = "Hayır"
$ws4.Range("C182").Value = "Motorlu Giriş Ünitesi"

$ws4.Range("D173").Value = 1
$ws4.Range("D174").Value = 2
$ws4.Range("D175").Value = 3
$ws4.Range("D176").Value = 4
$ws4.Range("D177").Value = 5
$ws4.Range("D178").Value = 6
$ws4.Range("D179").Value = 1
$ws4.Range("D180").Value = 2
$ws4.Range("D181").Value = 1
$ws4.Range("D182").Value = 2

# --- view state: selection moves to the new rows, active tab becomes "options" ---
$ws3.Range("Q62").Select()
$ws4.Activate()
$ws4.Range("A179").Select()
